$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(203957296, Omri Ben Shabat: -7,-5)"
$ws.Range("B1").Value = "(206532695, Matan Vakrat: 1,-7)"
$ws.Range("C1").Value = "(302962915, Asher  Odeh: -7,6)"
$ws.Range("D1").Value = "(308035542, Anastasia  Kubi: -7,-8)"
$ws.Range("E1").Value = "(311177802, Christina  Uksusman: 0,-4)"
$ws.Range("F1").Value = "(305251175, Or  Leder: -6,3)"

$ws.Range("A3").Value = "cost: 710.7430643061329"
$ws.Range("A4").Value = "time: 85.7178830382666"
